$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Style = "Normal"

$metaRange = $metaPara.Range
$metaStart = $metaRange.Start
$boldText = "Meta description"
$restText = ": Explore the secrets of ancient Egypt in Book of Anunnaki. Read our expert review, play for free, and learn about the game's theme, bonus features, payout, and more."
$metaRange.Text = $boldText + $restText

$boldRange = $d.Range($metaStart, $metaStart + $boldText.Length)
$boldRange.Bold = 1

# ------------------------------------------------------------------
# 2. Remove the duplicated bold "Play Book of Anunnaki..." paragraph
#    that used to sit near the end of the document, right before the
#    italic meta-description paragraph (now the last two paragraphs).
# ------------------------------------------------------------------
$paraCount = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs.Item($paraCount - 1)
$dupTitlePara.Range.Delete()

# ------------------------------------------------------------------
# 3. Replace the italic meta-description text at the very end of the
#    document with the new image-brief text, keeping its formatting.
#    (Scoped to the last paragraph's own Range so the earlier "Meta
#    description" paragraph - which contains the same old sentence -
#    is left untouched.)
# ------------------------------------------------------------------
$newDesc = "Create a feature image that will catch the attention of players of Book of Anunnaki. The image should be in cartoon style and have a happy Maya warrior wearing glasses. The warrior should be standing in front of an ancient Egyptian temple, holding the Book of Anunnaki in one hand while smiling at the rewards in the other hand. The background should be a desert scene with pyramids and camels. Make sure to include the game logo and name in the image to make it recognizable to players. The image should be colorful and playful to attract potential players to try the game."

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastParaRange = $lastPara.Range
$lastParaRange.MoveEnd(1, -1) | Out-Null
$lastParaRange.Find.Execute("Explore the secrets of ancient Egypt in Book of Anunnaki. Read our expert review, play for free, and learn about the game's theme, bonus features, payout, and more.", $true, $false, $false, $false, $false, $true, 1, $false, $newDesc, 2) | Out-Null

Write-Output "done"
